$wb = $excel.ActiveWorkbook

# --- Sheet "لها" (sheet1): remove the "مارك & وفر" credit line (was row 7, value 220). ---
# It moves to the "عليها" sheet instead (see below), so delete it here; everything below
# shifts up one row (old row 8 "178" becomes row 7, etc., total moves from row17 -> row16).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(7).Delete()

# Make room below the existing data (now ending at row 12) for six new May/June budget
# lines, while keeping a one-row gap both above (row 13) and below (row 20) them, same as
# the blank-row spacing the sheet already had before the edit. Inserting 5 rows at 16:20
# pushes the current total (currently on row 16) down to row 21, and leaves rows 14-19 free.
$ws1.Range("A16:A20").EntireRow.Insert()

$ws1.Range("A14").Value = "مصاريف أيسل (مايو)"
$ws1.Range("B14").Value = 7000
$ws1.Range("A15").Value = "قسط السيارة (مايو)"
$ws1.Range("B15").Value = 1000
$ws1.Range("A16").Value = "قسط الجمعية (مايو)"
$ws1.Range("B16").Value = 1000
$ws1.Range("A17").Value = "مصاريف أيسل (يونيو)"
$ws1.Range("B17").Value = 7000
$ws1.Range("A18").Value = "قسط السيارة (يونيو)"
$ws1.Range("B18").Value = 1000
$ws1.Range("A19").Value = "قسط الجمعية (يونيو)"
$ws1.Range("B19").Value = 1000

# --- Sheet "عليها" (sheet2): add back the "مارك & وفر" credit line as a new row 3. ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = "مشتريات بالكريدت (مارك & وفر)"
$ws2.Range("B3").Value = 220

# Update the manually-entered total on "عليها" (old 27700 + the 220 moved in above).
$ws2.Range("B8").Value = 27920

# --- Selections / active sheet, matching the saved view state. ---
$ws1.Activate()
$ws1.Range("F18").Select()

$ws2.Activate()
$ws2.Range("E10").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("C15").Select()
